$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns (zh-cn, de-de) and the generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 18:43:43"

# zh-cn sheet: status + Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 18:43:38"

# de-de sheet: status (Latest Handoff Datetime unchanged here)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 18:43:43"

# Columns holding the changed (now longer) status text need to widen to
# fit the new content - autofit, then settle on the fitted width.
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsZhCn.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(3).AutoFit()

$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38
